# Applies the pseudocode_flowchart.docx "chuyen doi tien te" edit:
#   1. Inserts a title block ("Thuat toan chuyen doi tien te") before the
#      existing flowchart content.
#   2. Rewrites the pseudo-code text run-by-run:
#        Input a,b,c        -> Input / dola   (2 runs)
#        Tb= (a+b+c)/3      -> vnd=dola*23000 (w:noProof rPr added)
#        Display Tb         -> Display / vnd  (2 runs)
#   3. Moves the stray "_GoBack" bookmark from after "end" to right after
#      the new "vnd" run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: replace the contents of $range (an existing Range) with the
# supplied run-level WordprocessingML, via a flat-OPC InsertXML payload.
# (NB: always build the payload into a plain variable first, then pass
#  only bare variables positionally -- passing a parenthesized/ concatenated
#  expression as a non-first argument after a COM-object argument is
#  mis-parsed by this interpreter.)
# ---------------------------------------------------------------------
function Set-RangeRunXml {
    param($range, [string]$runXml)

    $opcHead = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>'
    $opcHead = $opcHead + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
    $opcHead = $opcHead + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">'
    $opcHead = $opcHead + '<pkg:xmlData>'
    $opcHead = $opcHead + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
    $opcHead = $opcHead + '<w:body><w:p>'

    $opcTail = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $flatOpc = $opcHead + $runXml + $opcTail

    $range.InsertXML($flatOpc)
}

# ---------------------------------------------------------------------
# 0) The "_GoBack" bookmark currently sits after "end"; drop it here so
#    it can be re-created in its new home below (same w:id="0").
# ---------------------------------------------------------------------
$hasGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hasGoBack) {
    $goBackBm = $d.Bookmarks.Item("_GoBack")
    $goBackBm.Delete()
}

# ---------------------------------------------------------------------
# 1) Insert the new title block at the very start of the document:
#      <w:p/>
#      <w:p><w:pPr><sz 28/></w:pPr><w:r><sz 28/>Thuat toan...</w:r></w:p>
#      <w:p/>
# ---------------------------------------------------------------------
$titleText = "Thu" + [char]0x1EAD + "t to" + [char]0x00E1 + "n chuy" + [char]0x1EC3 + "n "
$titleText = $titleText + [char]0x0111 + [char]0x1ED5 + "i ti" + [char]0x1EC1 + "n t" + [char]0x1EC7

$titleBodyXml = '<w:p/>'
$titleBodyXml = $titleBodyXml + '<w:p><w:pPr><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr>'
$titleBodyXml = $titleBodyXml + '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>' + $titleText + '</w:t></w:r></w:p>'
$titleBodyXml = $titleBodyXml + '<w:p/>'

$titleOpcHead = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>'
$titleOpcHead = $titleOpcHead + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
$titleOpcHead = $titleOpcHead + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">'
$titleOpcHead = $titleOpcHead + '<pkg:xmlData>'
$titleOpcHead = $titleOpcHead + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$titleOpcHead = $titleOpcHead + '<w:body>'

$titleOpcTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$titleFlatOpc = $titleOpcHead + $titleBodyXml + $titleOpcTail

$docStart = $d.Range(0, 0)
$docStart.InsertXML($titleFlatOpc)

# ---------------------------------------------------------------------
# 2) Locate the 5 "pseudocode" paragraphs (Begin / Input / Tb= / Display /
#    end) by their text and rewrite the 3 that changed. Using text lookup
#    (rather than a hard-coded index) keeps this robust to the 3 extra
#    paragraphs just inserted above.
# ---------------------------------------------------------------------
$paragraphs = $d.Paragraphs
$count = $paragraphs.Count

$inputPara   = $null
$tbPara      = $null
$displayPara = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $paragraphs.Item($i)
    $fullText = $p.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph mark (chr 13).
    $t = $fullText.TrimEnd([char]13, [char]7)
    if ($t -eq "Input a,b,c") { $inputPara = $p }
    elseif ($t -eq "Tb= (a+b+c)/3") { $tbPara = $p }
    elseif ($t -eq "Display Tb") { $displayPara = $p }
}

# -- "Input a,b,c" -> "Input " + "dola" (2 runs) --
$inputStart = $inputPara.Range.Start
$inputEnd = $inputPara.Range.End - 1
$inputRange = $d.Range($inputStart, $inputEnd)
$inputRunXml = '<w:r><w:t xml:space="preserve">Input </w:t></w:r>'
$inputRunXml = $inputRunXml + '<w:r><w:t>dola</w:t></w:r>'
Set-RangeRunXml $inputRange $inputRunXml

# -- "Tb= (a+b+c)/3" -> "vnd=dola*23000" (single run, adds w:noProof) --
$tbStart = $tbPara.Range.Start
$tbEnd = $tbPara.Range.End - 1
$tbRange = $d.Range($tbStart, $tbEnd)
$tbRunXml = '<w:r><w:rPr><w:noProof/></w:rPr><w:t>vnd=dola*23000</w:t></w:r>'
Set-RangeRunXml $tbRange $tbRunXml

# -- "Display Tb" -> "Display " + "vnd" (2 runs) + relocated bookmark --
$displayStart = $displayPara.Range.Start
$displayEnd = $displayPara.Range.End - 1
$displayRange = $d.Range($displayStart, $displayEnd)
$displayRunXml = '<w:r><w:t xml:space="preserve">Display </w:t></w:r>'
$displayRunXml = $displayRunXml + '<w:r><w:t>vnd</w:t></w:r>'
$displayRunXml = $displayRunXml + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Set-RangeRunXml $displayRange $displayRunXml

Write-Output "done"
